{"js": "// The document starts with a centered \"<date> <weekday>\" title paragraph,\n// followed by a single 20-row x 5-column table of two-digit multiplication\n// problems (\"NN\u00d7NN=\"). The edit swaps the title's date/weekday and replaces\n// every problem in the table (each cell keeps its own run formatting).\n\n// 1) Update the date/weekday title (first paragraph of the body).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst dateOld = \"2023-07-01 Saturday\";\nconst dateNew = \"2023-07-02 Sunday\";\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text === dateOld) {\n    p.insertText(dateNew, Word.InsertLocation.replace);\n    break;\n  }\n}\n\n// 2) Replace every multiplication problem in the table, row by row.\n// Assigning `Table.values` rewrites each cell's text in place while\n// preserving the cell's existing paragraph/run formatting (font, size,\n// alignment, etc.), so only the `<w:t>` contents change.\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst newValues = [\n  [\"36\u00d768=\", \"30\u00d732=\", \"80\u00d763=\", \"19\u00d773=\", \"41\u00d713=\"],\n  [\"15\u00d741=\", \"40\u00d727=\", \"85\u00d750=\", \"34\u00d753=\", \"50\u00d743=\"],\n  [\"77\u00d775=\", \"58\u00d724=\", \"22\u00d710=\", \"49\u00d757=\", \"10\u00d718=\"],\n  [\"29\u00d743=\", \"73\u00d711=\", \"16\u00d752=\", \"67\u00d724=\", \"80\u00d779=\"],\n  [\"23\u00d732=\", \"76\u00d745=\", \"55\u00d769=\", \"97\u00d727=\", \"64\u00d782=\"],\n  [\"59\u00d784=\", \"72\u00d750=\", \"81\u00d775=\", \"16\u00d791=\", \"67\u00d732=\"],\n  [\"10\u00d756=\", \"66\u00d736=\", \"40\u00d7100=\", \"91\u00d740=\", \"29\u00d749=\"],\n  [\"63\u00d788=\", \"86\u00d792=\", \"66\u00d769=\", \"77\u00d797=\", \"69\u00d715=\"],\n  [\"29\u00d787=\", \"29\u00d719=\", \"55\u00d711=\", \"15\u00d766=\", \"73\u00d712=\"],\n  [\"25\u00d719=\", \"28\u00d776=\", \"49\u00d756=\", \"72\u00d720=\", \"11\u00d776=\"],\n  [\"61\u00d791=\", \"75\u00d776=\", \"81\u00d720=\", \"91\u00d754=\", \"60\u00d734=\"],\n  [\"25\u00d736=\", \"68\u00d746=\", \"46\u00d796=\", \"44\u00d761=\", \"94\u00d737=\"],\n  [\"38\u00d728=\", \"82\u00d772=\", \"34\u00d750=\", \"60\u00d737=\", \"89\u00d743=\"],\n  [\"70\u00d779=\", \"48\u00d768=\", \"69\u00d791=\", \"28\u00d766=\", \"16\u00d787=\"],\n  [\"68\u00d780=\", \"89\u00d716=\", \"76\u00d793=\", \"22\u00d767=\", \"89\u00d799=\"],\n  [\"21\u00d735=\", \"43\u00d745=\", \"22\u00d778=\", \"39\u00d722=\", \"22\u00d783=\"],\n  [\"33\u00d748=\", \"84\u00d752=\", \"64\u00d734=\", \"53\u00d794=\", \"72\u00d739=\"],\n  [\"37\u00d715=\", \"81\u00d757=\", \"92\u00d734=\", \"29\u00d766=\", \"98\u00d734=\"],\n  [\"39\u00d713=\", \"20\u00d767=\", \"45\u00d793=\", \"38\u00d774=\", \"24\u00d798=\"],\n  [\"69\u00d756=\", \"87\u00d710=\", \"91\u00d787=\", \"46\u00d739=\", \"71\u00d773=\"]\n];\n\nconst table = tables.items[0];\ntable.values = newValues;\n\nawait context.sync();\n", "ps1": "# The document starts with a centered \"<date> <weekday>\" title paragraph,\n# followed by a single 20-row x 5-column table of two-digit multiplication\n# problems (\"NN\u00d7NN=\"). The edit swaps the title's date/weekday and replaces\n# every problem in the table (each cell keeps its own run formatting).\n$d = $word.ActiveDocument\n\n# 1) Update the date/weekday title. Find the title text and replace just\n# that range so the paragraph's run formatting (font/size) is untouched.\n$titleRange = $d.Content\n$find = $titleRange.Find\n$find.Text = \"2023-07-01 Saturday\"\nif ($find.Execute()) {\n    $titleRange.Text = \"2023-07-02 Sunday\"\n}\n\n# 2) Replace every multiplication problem in the table, row by row.\n# Setting Cell.Range.Text rewrites only the cell's <w:t> contents and\n# preserves the existing paragraph/run formatting (font, size, alignment).\n$t = $d.Tables.Item(1)\n\n$newValues = @(\n    @(\"36\u00d768=\", \"30\u00d732=\", \"80\u00d763=\", \"19\u00d773=\", \"41\u00d713=\"),\n    @(\"15\u00d741=\", \"40\u00d727=\", \"85\u00d750=\", \"34\u00d753=\", \"50\u00d743=\"),\n    @(\"77\u00d775=\", \"58\u00d724=\", \"22\u00d710=\", \"49\u00d757=\", \"10\u00d718=\"),\n    @(\"29\u00d743=\", \"73\u00d711=\", \"16\u00d752=\", \"67\u00d724=\", \"80\u00d779=\"),\n    @(\"23\u00d732=\", \"76\u00d745=\", \"55\u00d769=\", \"97\u00d727=\", \"64\u00d782=\"),\n    @(\"59\u00d784=\", \"72\u00d750=\", \"81\u00d775=\", \"16\u00d791=\", \"67\u00d732=\"),\n    @(\"10\u00d756=\", \"66\u00d736=\", \"40\u00d7100=\", \"91\u00d740=\", \"29\u00d749=\"),\n    @(\"63\u00d788=\", \"86\u00d792=\", \"66\u00d769=\", \"77\u00d797=\", \"69\u00d715=\"),\n    @(\"29\u00d787=\", \"29\u00d719=\", \"55\u00d711=\", \"15\u00d766=\", \"73\u00d712=\"),\n    @(\"25\u00d719=\", \"28\u00d776=\", \"49\u00d756=\", \"72\u00d720=\", \"11\u00d776=\"),\n    @(\"61\u00d791=\", \"75\u00d776=\", \"81\u00d720=\", \"91\u00d754=\", \"60\u00d734=\"),\n    @(\"25\u00d736=\", \"68\u00d746=\", \"46\u00d796=\", \"44\u00d761=\", \"94\u00d737=\"),\n    @(\"38\u00d728=\", \"82\u00d772=\", \"34\u00d750=\", \"60\u00d737=\", \"89\u00d743=\"),\n    @(\"70\u00d779=\", \"48\u00d768=\", \"69\u00d791=\", \"28\u00d766=\", \"16\u00d787=\"),\n    @(\"68\u00d780=\", \"89\u00d716=\", \"76\u00d793=\", \"22\u00d767=\", \"89\u00d799=\"),\n    @(\"21\u00d735=\", \"43\u00d745=\", \"22\u00d778=\", \"39\u00d722=\", \"22\u00d783=\"),\n    @(\"33\u00d748=\", \"84\u00d752=\", \"64\u00d734=\", \"53\u00d794=\", \"72\u00d739=\"),\n    @(\"37\u00d715=\", \"81\u00d757=\", \"92\u00d734=\", \"29\u00d766=\", \"98\u00d734=\"),\n    @(\"39\u00d713=\", \"20\u00d767=\", \"45\u00d793=\", \"38\u00d774=\", \"24\u00d798=\"),\n    @(\"69\u00d756=\", \"87\u00d710=\", \"91\u00d787=\", \"46\u00d739=\", \"71\u00d773=\")\n)\n\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    for ($c = 1; $c -le $t.Columns.Count; $c++) {\n        $t.Cell($r, $c).Range.Text = $newValues[$r - 1][$c - 1]\n    }\n}\n\n"}
